$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "time_taken" timestamps (column F) on the existing "data"
#    sheet for rows 2-11.
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value = "2021-10-05 14:35:41.613043"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:41.613051"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:41.613054"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:41.613057"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:41.613059"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:41.613062"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:41.613065"
$dataSheet.Range("F9").Value = "2021-10-05 14:35:41.613067"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:41.613070"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:41.613073"

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" sheet right after "data" and populate it.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the bold/bordered/centered header style already used by the "data"
# sheet (style applied to B1 there) by copying it across the header row and
# the index cell, then overwriting the copied values.
$dataSheet.Range("B1").Copy($newSheet.Range("B1:G1"))
$dataSheet.Range("A2").Copy($newSheet.Range("A2"))

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Severe Combined Immunodeficiency (absent T absent B cells)"
$newSheet.Range("C2").Value = 234
# Keep data_version as the literal text "1.1" (not converted to a number).
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.1"
$newSheet.Range("E2").Value = "2021-08-26T07:23:02.567545Z"
$newSheet.Range("F2").Value = "2021-10-05 14:35:41.609354"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/234/?format=json"

# Make sure "data" stays the active tab, matching the original workbook.
$dataSheet.Activate()
